# Applies the "cryptos list" refresh: updated Price (D) / Volume(1h) (E)
# figures for most rows, plus a swap of the ImmutableX / BitcoinCash rows
# (23 and 24) which traded places in the new ranking.
#
# NumberFormat is forced to text ("@") before writing any Price cell whose
# new value looks like a plain number (e.g. "307.85", "1.00", "0.0800").
# Without this, Excel's automatic type detection would convert the text
# into a floating point number, which both loses significant trailing
# zeros (e.g. "1.00" -> "1") and introduces binary floating point noise
# (e.g. "0.831" -> "0.83099999999999996"). Forcing text keeps the values
# byte-identical to the source data, matching the original inline-string
# cell type.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "44.279.19"
$ws.Range("E2").Value = "  +2.49%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.428.70"
$ws.Range("E3").Value = "  +1.93%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "307.85"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "99.92"
$ws.Range("E6").Value = "  +2.96%  "
$ws.Range("E7").Value = "  +0.80%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.499"
$ws.Range("E9").Value = "  -0.35%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.38"
$ws.Range("E10").Value = "  +3.66%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0800"
$ws.Range("E11").Value = "  +1.28%  "
$ws.Range("E12").Value = "  +2.55%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.74"
$ws.Range("E13").Value = "  +1.58%  "
$ws.Range("E14").Value = "  +2.19%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.807.89"
$ws.Range("E15").Value = "  +2.05%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.450.56"
$ws.Range("E16").Value = "  +2.88%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.831"
$ws.Range("E17").Value = "  +2.79%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "44.246.94"
$ws.Range("E18").Value = "  +2.49%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.32"
$ws.Range("E19").Value = "  +0.50%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.45"
$ws.Range("E20").Value = "  +1.59%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0907"
$ws.Range("E21").Value = "  +1.87%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.70"
$ws.Range("E22").Value = "  +0.21%  "
$ws.Range("B23").Value = "ImmutableX"
$ws.Range("C23").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.31"
$ws.Range("E23").Value = "  +3.61%  "
$ws.Range("B24").Value = "BitcoinCash"
$ws.Range("C24").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "240.38"
$ws.Range("E24").Value = "  +2.07%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.47"
$ws.Range("E25").Value = "  +1.24%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.29"
$ws.Range("E27").Value = "  +1.86%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.32"
$ws.Range("E28").Value = "  -1.98%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.50"
$ws.Range("E29").Value = "  +4.03%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "32.92"
$ws.Range("E30").Value = "  +4.20%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.119"
$ws.Range("E31").Value = "  +15.69%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "18.67"
$ws.Range("E32").Value = "  +8.60%  "
$ws.Range("E33").Value = "  +1.62%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.00"
$ws.Range("E34").Value = "  +0.07%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0763"
$ws.Range("E35").Value = "  +3.59%  "
$ws.Range("E36").Value = "  +3.36%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.56"
$ws.Range("E37").Value = "  +5.65%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "129.67"
$ws.Range("E38").Value = "  +22.03%  "
$ws.Range("E39").Value = "  +4.59%  "
$ws.Range("E40").Value = "  -0.75%  "
$ws.Range("E41").Value = "  +0.04%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "21.17"
$ws.Range("E42").Value = "  -5.52%  "
$ws.Range("E43").Value = "  +2.74%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.957.63"
$ws.Range("E44").Value = "  +0.00%  "
$ws.Range("E45").Value = "  +1.96%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.88"
$ws.Range("E46").Value = "  +4.56%  "
$ws.Range("E47").Value = "  +1.03%  "
$ws.Range("E48").Value = "  +8.87%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.670.31"
$ws.Range("E49").Value = "  +2.07%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "53.48"
$ws.Range("E50").Value = "  +1.22%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "73.74"
$ws.Range("E51").Value = "  +2.33%  "
